$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild row 15 as a clean, unformatted row (matching row 14's pattern):
# drop the old placeholder row's per-cell styles/row attributes by
# deleting the row and inserting a fresh blank one in its place.
$ws.Rows("15:15").Delete()
$ws.Rows("15:15").Insert()

$ws.Range("B15").Value = "Jhan Carlos Ortiz"
$ws.Range("C15").Value = "Tarjeta de identidad"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64646464"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "15"
$ws.Range("F15").Value = "Villa del lago"
$ws.Range("G15").Value = "call 12 12 12"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "646464"
